# modified the BorderTemplate class.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 - Read Value / ValueInput: Read() -> Display(), QuickRead() -> QuickDisplay()
$ws.Range("D11").Value = "Display()"
$ws.Range("F11").Value = "QuickDisplay()"

# Row 12 - Read List of Values / ListInput: Read() -> Display(), QuickRead() -> QuickDisplay()
$ws.Range("D12").Value = "Display()"
$ws.Range("F12").Value = "QuickDisplay()"

# Row 13 - Display Value / ValueOutput: Write<T>() -> Display(), QuickWrite<T>() -> QuickDisplay()
$ws.Range("D13").Value = "Display()"
$ws.Range("F13").Value = "QuickDisplay()"

# Row 14 - Display List of Values: Write<T>() -> Display(), QuickWrite<T>() -> QuickDisplay()
$ws.Range("D14").Value = "Display()"
$ws.Range("F14").Value = "QuickDisplay()"

# Row 15 - Yes/No Question: ReadAnswer() -> Display() / ReadAnswer()
$ws.Range("D15").Value = "Display() / ReadAnswer()"

# Row 14 - OutputList -> ListOutput (new shared string added after "Display() / ReadAnswer()")
$ws.Range("C14").Value = "ListOutput"

# Row 16 - Pause Control: QuickPause() -> QuickDisplay()
$ws.Range("F16").Value = "QuickDisplay()"

# Update the active cell/selection in the sheet view (cosmetic, matches diff)
$ws.Range("F17").Select()
